$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2 and 3 (refresh timestamps/durations)
$ws.Range("C2").Value = "2024-09-03 14:24:16"
$ws.Range("K2").Value = "23h 59m 43s"

$ws.Range("C3").Value = "2024-09-03 14:24:16"
$ws.Range("K3").Value = "23h 59m 38s"

# Add new rows 4 and 5
$ws.Range("A4").Value = "osano_consentmanager"
$ws.Range("B4").Value = ".victorinsurance.nl"
$ws.Range("C4").Value = "2024-09-03 14:24:48"
$ws.Range("D4").Value = $true
$ws.Range("E4").Value = "Yes"
$ws.Range("F4").Value = "Yes"
$ws.Range("G4").Value = "Osano"
$ws.Range("H4").Value = "Yes"
$ws.Range("I4").Value = "Type1 (Manage Cookies)"
$ws.Range("J4").Value = "No"
$ws.Range("K4").Value = "23h 59m 44s"

$ws.Range("A5").Value = "osano_consentmanager_uuid"
$ws.Range("B5").Value = ".victorinsurance.nl"
$ws.Range("C5").Value = "2024-09-03 14:24:48"
$ws.Range("D5").Value = $true
$ws.Range("E5").Value = "Yes"
$ws.Range("F5").Value = "Yes"
$ws.Range("G5").Value = "Osano"
$ws.Range("H5").Value = "Yes"
$ws.Range("I5").Value = "Type1 (Manage Cookies)"
$ws.Range("J5").Value = "No"
$ws.Range("K5").Value = "23h 59m 39s"

# Reset column widths to default (remove explicit <cols> customWidth entries)
$ws.Columns("A:K").ColumnWidth = 8
